$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '75.778.73'
$ws.Range('E2').Value = '  +9.12%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.736.58'
$ws.Range('E3').Value = '  +12.63%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '189.55'
$ws.Range('E5').Value = '  +13.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '594.36'
$ws.Range('E6').Value = '  +5.39%  '

$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.545'
$ws.Range('E8').Value = '  +6.15%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.199'
$ws.Range('E9').Value = '  +17.30%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.740.96'
$ws.Range('E10').Value = '  +12.85%  '

$ws.Range('E11').Value = '  +1.52%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.364'
$ws.Range('E12').Value = '  +9.08%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.80'
$ws.Range('E13').Value = '  +2.79%  '

$ws.Range('B14').Value = 'WrappedBTC'
$ws.Range('C14').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '75.593.82'
$ws.Range('E14').Value = '  +9.07%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.183.33'
$ws.Range('E15').Value = '  +10.73%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000192'
$ws.Range('E16').Value = '  +7.95%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.12'
$ws.Range('E17').Value = '  +13.38%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.710.85'
$ws.Range('E18').Value = '  +11.71%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.61'
$ws.Range('E19').Value = '  +34.57%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.27'
$ws.Range('E20').Value = '  +13.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '380.83'
$ws.Range('E21').Value = '  +10.98%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.33'
$ws.Range('E22').Value = '  +18.44%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.12'
$ws.Range('E23').Value = '  +6.72%  '

$ws.Range('E24').Value = '  +5.24%  '

$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '71.14'
$ws.Range('E25').Value = '  +7.68%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.27'
$ws.Range('E27').Value = '  +12.01%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.68'
$ws.Range('E28').Value = '  +14.18%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.830.12'
$ws.Range('E29').Value = '  +10.81%  '

$ws.Range('E30').Value = '  -3.58%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0996'
$ws.Range('E31').Value = '  +17.09%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '530.78'
$ws.Range('E32').Value = '  +16.80%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.43'
$ws.Range('E33').Value = '  +14.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.96'
$ws.Range('E34').Value = '  +7.69%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.79'
$ws.Range('E35').Value = '  +10.27%  '

$ws.Range('E36').Value = '  -0.06%  '

$ws.Range('E37').Value = '  +9.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.26'
$ws.Range('E38').Value = '  +1.96%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.61'
$ws.Range('E39').Value = '  +7.50%  '

$ws.Range('E40').Value = '  +1.53%  '

$ws.Range('E41').Value = '  +0.02%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.13'
$ws.Range('E42').Value = '  +16.28%  '

$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '174.40'
$ws.Range('E43').Value = '  +28.31%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.74'
$ws.Range('E44').Value = '  +14.33%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.336'
$ws.Range('E45').Value = '  +10.99%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.23'
$ws.Range('E46').Value = '  +14.07%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.46'
$ws.Range('E47').Value = '  +17.98%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '39.21'
$ws.Range('E48').Value = '  +3.64%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0854'
$ws.Range('E49').Value = '  +18.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.73'
$ws.Range('E50').Value = '  +10.00%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.552'
$ws.Range('E51').Value = '  +12.90%  '
